# Re-run of 380 kV case: update power-flow result values (MW) in pl_mw
# for data rows 0-23 (sheet rows 2-25), columns B,C,D,F,G,J,K,M,O.
# Columns A,E,H,I,L,N and the header row are unchanged.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.339877782745333
$ws.Range("C2").Value = 0.03220462881643016
$ws.Range("D2").Value = 0.2303704542396048
$ws.Range("F2").Value = 1.724232725430632
$ws.Range("G2").Value = 0.002484166947418612
$ws.Range("J2").Value = 0.309022512547422
$ws.Range("K2").Value = 0.2949386773609888
$ws.Range("M2").Value = 0.2695704619020916
$ws.Range("O2").Value = 4.149581285493582

# Row 3
$ws.Range("B3").Value = 0.308052815535433
$ws.Range("C3").Value = 0.02885977443804677
$ws.Range("D3").Value = 0.2255462642936834
$ws.Range("F3").Value = 1.728584257974674
$ws.Range("G3").Value = 0.002486739840041292
$ws.Range("J3").Value = 0.3073761614179844
$ws.Range("K3").Value = 0.2623871385593333
$ws.Range("M3").Value = 0.2568595057164842
$ws.Range("O3").Value = 4.171451726982724

# Row 4
$ws.Range("B4").Value = 0.2885786495893115
$ws.Range("C4").Value = 0.02679389792907472
$ws.Range("D4").Value = 0.2226794944728141
$ws.Range("F4").Value = 1.732099630468639
$ws.Range("G4").Value = 0.002488404259369859
$ws.Range("J4").Value = 0.3065196435658066
$ws.Range("K4").Value = 0.242413613616506
$ws.Range("M4").Value = 0.2491664214036362
$ws.Range("O4").Value = 4.186931279578786

# Row 5
$ws.Range("B5").Value = 0.280659925413687
$ws.Range("C5").Value = 0.02594903660660464
$ws.Range("D5").Value = 0.2215353315272353
$ws.Range("F5").Value = 1.73374438013262
$ws.Range("G5").Value = 0.002489103873530348
$ws.Range("J5").Value = 0.3062094421178188
$ws.Range("K5").Value = 0.2342780105601605
$ws.Range("M5").Value = 0.2460596383375844
$ws.Range("O5").Value = 4.193755105723639

# Row 6
$ws.Range("B6").Value = 0.2793460764014242
$ws.Range("C6").Value = 0.02580856827569278
$ws.Range("D6").Value = 0.2213468009315847
$ws.Range("F6").Value = 1.73403030923874
$ws.Range("G6").Value = 0.002489221335188177
$ws.Range("J6").Value = 0.3061602797574068
$ws.Range("K6").Value = 0.2329273405536298
$ws.Range("M6").Value = 0.2455454684367453
$ws.Range("O6").Value = 4.194919350047783

# Row 7
$ws.Range("B7").Value = 0.2884717848591265
$ws.Range("C7").Value = 0.02678251591699166
$ws.Range("D7").Value = 0.2226639662978158
$ws.Range("F7").Value = 1.732120952792009
$ws.Range("G7").Value = 0.002488413608064693
$ws.Range("J7").Value = 0.3065153028002996
$ws.Range("K7").Value = 0.2423038781791576
$ws.Range("M7").Value = 0.2491244077814301
$ws.Range("O7").Value = 4.187021219655151

# Row 8
$ws.Range("B8").Value = 0.3288909847579191
$ws.Range("C8").Value = 0.03105386698635471
$ws.Range("D8").Value = 0.228687354216035
$ws.Range("F8").Value = 1.725558118520027
$ws.Range("G8").Value = 0.002485036549072427
$ws.Range("J8").Value = 0.3084228290037387
$ws.Range("K8").Value = 0.2837124087326117
$ws.Range("M8").Value = 0.2651646840799771
$ws.Range("O8").Value = 4.156696545142836

# Row 9
$ws.Range("B9").Value = 0.4086654251378832
$ws.Range("C9").Value = 0.03933205527995653
$ws.Range("D9").Value = 0.2412516033652281
$ws.Range("F9").Value = 1.719378503299083
$ws.Range("G9").Value = 0.002479082942706391
$ws.Range("J9").Value = 0.3133876729462344
$ws.Range("K9").Value = 0.3650043777568044
$ws.Range("M9").Value = 0.2974984095796245
$ws.Range("O9").Value = 4.113504442957947

# Row 10
$ws.Range("B10").Value = 0.4675739282299958
$ws.Range("C10").Value = 0.04535259879959597
$ws.Range("D10").Value = 0.2509373657687632
$ws.Range("F10").Value = 1.718915229736609
$ws.Range("G10").Value = 0.002475112495675511
$ws.Range("J10").Value = 0.3177818330127309
$ws.Range("K10").Value = 0.4247700482792425
$ws.Range("M10").Value = 0.3217845540137318
$ws.Range("O10").Value = 4.09169688445607

# Row 11
$ws.Range("B11").Value = 0.4944350607725312
$ws.Range("C11").Value = 0.04807782827298013
$ws.Range("D11").Value = 0.2554416812344442
$ws.Range("F11").Value = 1.71958956267072
$ws.Range("G11").Value = 0.002473393027996912
$ws.Range("J11").Value = 0.3199430398515233
$ws.Range("K11").Value = 0.4519651265108848
$ws.Range("M11").Value = 0.332947227839874
$ws.Range("O11").Value = 4.083932532823212

# Row 12
$ws.Range("B12").Value = 0.5046154284423778
$ws.Range("C12").Value = 0.04910781367222228
$ws.Range("D12").Value = 0.2571613828896773
$ws.Range("F12").Value = 1.719972135816207
$ws.Range("G12").Value = 0.002472754313770475
$ws.Range("J12").Value = 0.3207847588999329
$ws.Range("K12").Value = 0.4622638987468122
$ws.Range("M12").Value = 0.3371906098960622
$ws.Range("O12").Value = 4.081302462468955

# Row 13
$ws.Range("B13").Value = 0.5024225273463969
$ws.Range("C13").Value = 0.04888607765535369
$ws.Range("D13").Value = 0.256790392758333
$ws.Range("F13").Value = 1.719884084679151
$ws.Range("G13").Value = 0.002472891321087945
$ws.Range("J13").Value = 0.3206024429194798
$ws.Range("K13").Value = 0.4600458553298381
$ws.Range("M13").Value = 0.3362759987835844
$ws.Range("O13").Value = 4.08185510105227

# Row 14
$ws.Range("B14").Value = 0.4952724352880296
$ws.Range("C14").Value = 0.04816260616784973
$ws.Range("D14").Value = 0.2555828818577908
$ws.Range("F14").Value = 1.719618487875962
$ws.Range("G14").Value = 0.002473340231853266
$ws.Range("J14").Value = 0.3200118213871832
$ws.Range("K14").Value = 0.4528124036594363
$ws.Range("M14").Value = 0.3332960074291051
$ws.Range("O14").Value = 4.083709939558304

# Row 15
$ws.Range("B15").Value = 0.4908939086374744
$ws.Range("C15").Value = 0.04771919701734362
$ws.Range("D15").Value = 0.2548450685509636
$ws.Range("F15").Value = 1.719472368041522
$ws.Range("G15").Value = 0.002473616817428
$ws.Range("J15").Value = 0.3196530846317955
$ws.Range("K15").Value = 0.4483817679657989
$ws.Range("M15").Value = 0.3314727964695194
$ws.Range("O15").Value = 4.084886471241532

# Row 16
$ws.Range("B16").Value = 0.4658197246775444
$ws.Range("C16").Value = 0.04517422161373474
$ws.Range("D16").Value = 0.2506449652464369
$ws.Range("F16").Value = 1.718888965816973
$ws.Range("G16").Value = 0.002475226606844657
$ws.Range("J16").Value = 0.317643856963457
$ws.Range("K16").Value = 0.4229928937823786
$ws.Range("M16").Value = 0.3210573408708015
$ws.Range("O16").Value = 4.092247695238342

# Row 17
$ws.Range("B17").Value = 0.4504534169789451
$ws.Range("C17").Value = 0.04360945508801706
$ws.Range("D17").Value = 0.2480934159931536
$ws.Range("F17").Value = 1.718757718061667
$ws.Range("G17").Value = 0.002476236327359179
$ws.Range("J17").Value = 0.316452811590139
$ws.Range("K17").Value = 0.4074191926971196
$ws.Range("M17").Value = 0.3146970628242443
$ws.Range("O17").Value = 4.097315838413095

# Row 18
$ws.Range("B18").Value = 0.4416211149366802
$ws.Range("C18").Value = 0.04270817083991574
$ws.Range("D18").Value = 0.2466350831160895
$ws.Range("F18").Value = 1.718765538919186
$ws.Range("G18").Value = 0.002476825256943378
$ws.Range("J18").Value = 0.3157830296286761
$ws.Range("K18").Value = 0.3984623341272879
$ws.Range("M18").Value = 0.3110496143289581
$ws.Range("O18").Value = 4.100433829472848

# Row 19
$ws.Range("B19").Value = 0.4386316926920415
$ws.Range("C19").Value = 0.04240279461397733
$ws.Range("D19").Value = 0.2461429087770597
$ws.Range("F19").Value = 1.718782497369389
$ws.Range("G19").Value = 0.002477026062292279
$ws.Range("J19").Value = 0.3155588771786455
$ws.Range("K19").Value = 0.3954298342666789
$ws.Range("M19").Value = 0.3098165141408984
$ws.Range("O19").Value = 4.101524378251298

# Row 20
$ws.Range("B20").Value = 0.4520885699641326
$ws.Range("C20").Value = 0.04377615917421451
$ws.Range("D20").Value = 0.2483640761288513
$ws.Range("F20").Value = 1.718763067620444
$ws.Range("G20").Value = 0.002476127996310973
$ws.Range("J20").Value = 0.3165780195493113
$ws.Range("K20").Value = 0.409076968441866
$ws.Range("M20").Value = 0.3153730072962517
$ws.Range("O20").Value = 4.096755322863203

# Row 21
$ws.Range("B21").Value = 0.4973723599146354
$ws.Range("C21").Value = 0.04837516196911906
$ws.Range("D21").Value = 0.255937177571596
$ws.Range("F21").Value = 1.719693047880924
$ws.Range("G21").Value = 0.002473208039601354
$ws.Range("J21").Value = 0.3201846685851564
$ws.Range("K21").Value = 0.4549370316132695
$ws.Range("M21").Value = 0.3341708616955401
$ws.Range("O21").Value = 4.083156711225655

# Row 22
$ws.Range("B22").Value = 0.5270180758039658
$ws.Range("C22").Value = 0.05136919561479658
$ws.Range("D22").Value = 0.2609682781682778
$ws.Range("F22").Value = 1.721042341034433
$ws.Range("G22").Value = 0.002471371993675229
$ws.Range("J22").Value = 0.3226777222059667
$ws.Range("K22").Value = 0.4849124328197547
$ws.Range("M22").Value = 0.3465513839245773
$ws.Range("O22").Value = 4.076076896261611

# Row 23
$ws.Range("B23").Value = 0.5111911666663502
$ws.Range("C23").Value = 0.04977230934748889
$ws.Range("D23").Value = 0.2582756505417336
$ws.Range("F23").Value = 1.720254369223127
$ws.Range("G23").Value = 0.00247234532795189
$ws.Range("J23").Value = 0.321334704005821
$ws.Range("K23").Value = 0.4689138605355936
$ws.Range("M23").Value = 0.3399350317034191
$ws.Range("O23").Value = 4.079690096261402

# Row 24
$ws.Range("B24").Value = 0.4513493110570437
$ws.Range("C24").Value = 0.04370079748329658
$ws.Range("D24").Value = 0.2482416839274748
$ws.Range("F24").Value = 1.718760389675907
$ws.Range("G24").Value = 0.002476176946593726
$ws.Range("J24").Value = 0.3165213664130846
$ws.Range("K24").Value = 0.4083274983667877
$ws.Range("M24").Value = 0.3150673842549381
$ws.Range("O24").Value = 4.097008095720241

# Row 25
$ws.Range("B25").Value = 0.3870309338013271
$ws.Range("C25").Value = 0.03710325106163737
$ws.Range("D25").Value = 0.2377724924019162
$ws.Range("F25").Value = 1.720334194409617
$ws.Range("G25").Value = 0.002480622371999392
$ws.Range("J25").Value = 0.3119134730981941
$ws.Range("K25").Value = 0.3430046099381343
$ws.Range("M25").Value = 0.2886577361006886
$ws.Range("O25").Value = 4.123446308897456
